# Update the "Förändrad" (Changed/Updated) date column C for all data rows
# from 45180 (2023-09-11) to 45181 (2023-09-12).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("C2:C532")
$rng.Value2 = 45181
